$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("questionId", "prompt", "answers", "topic", "subtopic", "difficulty")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("F1").Select() | Out-Null
